# Scheduled-runner market-data refresh: rewrite the static
# currentAveragePrice / LevePrice / LeveProfit columns (H:N) for the
# leves whose underlying item prices moved, across all 8 job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 2726.4285
$ws.Range("I5").Value = 2726.4285
$ws.Range("K5").Value = 2726.4285
$ws.Range("M5").Value = -2611.4285

$ws.Range("H6").Value = 81169.39999999999
$ws.Range("I6").Value = 81169.39999999999
$ws.Range("K6").Value = 243508.2
$ws.Range("M6").Value = -243396.2

$ws.Range("H17").Value = 3584.96
$ws.Range("J17").Value = 4326.2
$ws.Range("L17").Value = 12978.6
$ws.Range("N17").Value = -13314.6

$ws.Range("H33").Value = 802.3200000000001
$ws.Range("I33").Value = 537
$ws.Range("K33").Value = 537
$ws.Range("M33").Value = -308

$ws.Range("H62").Value = 8027.1943
$ws.Range("I62").Value = 8620.606
$ws.Range("J62").Value = 1499.6666
$ws.Range("K62").Value = 8620.606
$ws.Range("L62").Value = 1499.6666
$ws.Range("M62").Value = -7996.606
$ws.Range("N62").Value = -2747.6666

$ws.Range("H65").Value = 8027.1943
$ws.Range("I65").Value = 8620.606
$ws.Range("J65").Value = 1499.6666
$ws.Range("K65").Value = 43103.03
$ws.Range("L65").Value = 7498.333000000001
$ws.Range("M65").Value = -39983.03
$ws.Range("N65").Value = -13738.333

$ws.Range("H69").Value = 15060.929
$ws.Range("I69").Value = 9830.5
$ws.Range("J69").Value = 16487.408
$ws.Range("K69").Value = 29491.5
$ws.Range("L69").Value = 49462.224
$ws.Range("M69").Value = -28617.5
$ws.Range("N69").Value = -51210.224

$ws.Range("H72").Value = 15060.929
$ws.Range("I72").Value = 9830.5
$ws.Range("J72").Value = 16487.408
$ws.Range("K72").Value = 88474.5
$ws.Range("L72").Value = 148386.672
$ws.Range("M72").Value = -84106.5
$ws.Range("N72").Value = -157122.672

$ws.Range("H132").Value = 3284.4211
$ws.Range("I132").Value = 3284.4211
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9853.263300000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7323.263300000001
$ws.Range("N132").ClearContents()

$ws.Range("H137").Value = 4348.3516
$ws.Range("I137").Value = 3949.7334
$ws.Range("J137").Value = 6056.7144
$ws.Range("K137").Value = 11849.2002
$ws.Range("L137").Value = 18170.1432
$ws.Range("M137").Value = -9299.200199999999
$ws.Range("N137").Value = -23270.1432

$ws.Range("H138").Value = 7920.767
$ws.Range("J138").Value = 7877.2656
$ws.Range("L138").Value = 23631.7968
$ws.Range("N138").Value = -33911.7968

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2730.1628
$ws.Range("I2").Value = 2297.1714
$ws.Range("J2").Value = 4624.5
$ws.Range("K2").Value = 2297.1714
$ws.Range("L2").Value = 4624.5
$ws.Range("M2").Value = -2184.1714
$ws.Range("N2").Value = -4850.5

$ws.Range("H19").Value = 15501.5
$ws.Range("I19").Value = 15501.5
$ws.Range("K19").Value = 15501.5
$ws.Range("M19").Value = -15272.5

$ws.Range("H32").Value = 23803.637
$ws.Range("I32").Value = 23803.637
$ws.Range("K32").Value = 23803.637
$ws.Range("M32").Value = -23516.637

$ws.Range("H61").Value = 11747
$ws.Range("I61").Value = 7613.5
$ws.Range("J61").Value = 20014
$ws.Range("K61").Value = 7613.5
$ws.Range("L61").Value = 20014
$ws.Range("M61").Value = -7401.5
$ws.Range("N61").Value = -20438

$ws.Range("H88").Value = 1913.2858
$ws.Range("I88").Value = 1598.25
$ws.Range("J88").Value = 2333.3333
$ws.Range("K88").Value = 1598.25
$ws.Range("L88").Value = 2333.3333
$ws.Range("M88").Value = -1192.25
$ws.Range("N88").Value = -3145.3333

$ws.Range("H91").Value = 1913.2858
$ws.Range("I91").Value = 1598.25
$ws.Range("J91").Value = 2333.3333
$ws.Range("K91").Value = 1598.25
$ws.Range("L91").Value = 2333.3333
$ws.Range("M91").Value = -194.25
$ws.Range("N91").Value = -5141.3333

$ws.Range("H116").Value = 2730.1628
$ws.Range("I116").Value = 2297.1714
$ws.Range("J116").Value = 4624.5
$ws.Range("K116").Value = 2297.1714
$ws.Range("L116").Value = 4624.5
$ws.Range("M116").Value = -3.171400000000176
$ws.Range("N116").Value = -9212.5

$ws.Range("H132").Value = 4379.3
$ws.Range("I132").Value = 3336.775
$ws.Range("K132").Value = 10010.325
$ws.Range("M132").Value = -7480.325000000001

$ws.Range("H134").Value = 91640.664
$ws.Range("J134").Value = 91640.664
$ws.Range("L134").Value = 91640.664
$ws.Range("N134").Value = -101780.664

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H136").Value = 11747
$ws.Range("I136").Value = 7613.5
$ws.Range("J136").Value = 20014
$ws.Range("K136").Value = 22840.5
$ws.Range("L136").Value = 60042
$ws.Range("M136").Value = -20290.5
$ws.Range("N136").Value = -65142

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("M139").ClearContents()
$ws.Range("N139").ClearContents()

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2730.1628
$ws.Range("I3").Value = 2297.1714
$ws.Range("J3").Value = 4624.5
$ws.Range("K3").Value = 2297.1714
$ws.Range("L3").Value = 4624.5
$ws.Range("M3").Value = -2183.1714
$ws.Range("N3").Value = -4852.5

$ws.Range("H5").Value = 37305.332
$ws.Range("I5").Value = 250
$ws.Range("K5").Value = 250
$ws.Range("M5").Value = -137

$ws.Range("H86").Value = 11520.286
$ws.Range("I86").Value = 15550
$ws.Range("K86").Value = 15550
$ws.Range("M86").Value = -14427

$ws.Range("H89").Value = 11520.286
$ws.Range("I89").Value = 15550
$ws.Range("K89").Value = 77750
$ws.Range("M89").Value = -72134

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4709
$ws.Range("I58").Value = 4467.5454
$ws.Range("J58").Value = 5240.2
$ws.Range("K58").Value = 4467.5454
$ws.Range("L58").Value = 5240.2
$ws.Range("M58").Value = -4264.5454
$ws.Range("N58").Value = -5646.2

$ws.Range("H136").Value = 4709
$ws.Range("I136").Value = 4467.5454
$ws.Range("J136").Value = 5240.2
$ws.Range("K136").Value = 13402.6362
$ws.Range("L136").Value = 15720.6
$ws.Range("M136").Value = -10852.6362
$ws.Range("N136").Value = -20820.6

$ws.Range("H141").Value = 321247.47
$ws.Range("J141").Value = 399525.28
$ws.Range("L141").Value = 399525.28
$ws.Range("N141").Value = -409885.28

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1558.6316
$ws.Range("I5").Value = 1079.6
$ws.Range("J5").Value = 3355
$ws.Range("K5").Value = 3238.8
$ws.Range("L5").Value = 10065
$ws.Range("M5").Value = -3126.8
$ws.Range("N5").Value = -10289

$ws.Range("H6").Value = 75.5
$ws.Range("I6").Value = 75.5
$ws.Range("K6").Value = 226.5
$ws.Range("M6").Value = -113.5

$ws.Range("H128").Value = 265994.28
$ws.Range("I128").Value = 265994.28
$ws.Range("K128").Value = 797982.8400000001
$ws.Range("M128").Value = -793002.8400000001

$ws.Range("H135").Value = 1558.6316
$ws.Range("I135").Value = 1079.6
$ws.Range("J135").Value = 3355
$ws.Range("K135").Value = 9716.4
$ws.Range("L135").Value = 30195
$ws.Range("M135").Value = -7181.4
$ws.Range("N135").Value = -35265

$ws.Range("H140").Value = 2049.4146
$ws.Range("I140").Value = 1108.4667
$ws.Range("K140").Value = 3325.4001
$ws.Range("M140").Value = 1854.5999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3713.8572
$ws.Range("I80").Value = 3333
$ws.Range("J80").Value = 3999.5
$ws.Range("K80").Value = 3333
$ws.Range("L80").Value = 3999.5
$ws.Range("M80").Value = -2335
$ws.Range("N80").Value = -5995.5

$ws.Range("H83").Value = 3713.8572
$ws.Range("I83").Value = 3333
$ws.Range("J83").Value = 3999.5
$ws.Range("K83").Value = 16665
$ws.Range("L83").Value = 19997.5
$ws.Range("M83").Value = -11673
$ws.Range("N83").Value = -29981.5

$ws.Range("H138").Value = 208662.67
$ws.Range("J138").Value = 208662.67
$ws.Range("L138").Value = 208662.67
$ws.Range("N138").Value = -218942.67

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 14861.714
$ws.Range("I12").Value = 2084.3333
$ws.Range("J12").Value = 24444.75
$ws.Range("K12").Value = 2084.3333
$ws.Range("L12").Value = 24444.75
$ws.Range("M12").Value = -1914.3333
$ws.Range("N12").Value = -24784.75

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H122").Value = 4739.087
$ws.Range("I122").Value = 3999.95
$ws.Range("K122").Value = 11999.85
$ws.Range("M122").Value = -9549.849999999999

$ws.Range("H132").Value = 4688.543
$ws.Range("I132").Value = 4090.8
$ws.Range("K132").Value = 12272.4
$ws.Range("M132").Value = -9742.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws.Range("H126").Value = 2029.9318
$ws.Range("I126").Value = 1647.5
$ws.Range("J126").Value = 4452
$ws.Range("K126").Value = 4942.5
$ws.Range("L126").Value = 13356
$ws.Range("M126").Value = -2472.5
$ws.Range("N126").Value = -18296

$ws.Range("H132").Value = 6211.07
$ws.Range("I132").Value = 5953.2563
$ws.Range("K132").Value = 17859.7689
$ws.Range("M132").Value = -15329.7689

$ws.Range("H136").Value = 2668.434
$ws.Range("J136").Value = 6059
$ws.Range("L136").Value = 18177
$ws.Range("N136").Value = -23277

$ws.Range("H140").Value = 94059.71000000001
$ws.Range("J140").Value = 96403.5
$ws.Range("L140").Value = 96403.5
$ws.Range("N140").Value = -106763.5
